# "some changed rev 1.0"
#
# Updates two resistor-divider input values and lets the dependent formulas
# recalc on their own, then leaves the UI selection/active-tab the way the
# author left it when they saved.
#
#   Spannungsteiler!C2 : 82000 -> 5600   (C4, C5 are formulas and recalc)
#   6th tab  !B2        : 33000 -> 10000  (B4, B5 are formulas and recalc)
#   6th tab  !B3        : 17000 -> 5600
#
# Note: this workbook's 6th/7th tabs ("OpAmp" / "TLE42764") have their
# sheetId order swapped relative to tab order, so Worksheets.Item(<name>)
# resolves to the wrong underlying sheet data here. Item(<1-based tab
# position>) addresses the sheets reliably, so that's used for all of them.

$wb = $excel.ActiveWorkbook

# --- Spannungsteiler (5th tab) ---
$wsSpannung = $wb.Worksheets.Item(5)
$wsSpannung.Range("C2").Value = 5600
$wsSpannung.Range("K14").Select() | Out-Null

# --- 6th tab ---
$wsTle = $wb.Worksheets.Item(6)
$wsTle.Range("B2").Value = 10000
$wsTle.Range("B3").Value = 5600

# Select/activate the 6th tab last so it ends up the active sheet with its
# own selection, matching the saved state in the target workbook.
$wsTle.Range("B4").Select() | Out-Null
$wsTle.Activate()
